# Updated symbol list on Mon Dec 26 20:46:46 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as TEXT in this
# workbook (t="inlineStr"), so force text formatting before writing the
# new value - otherwise Excel would happily re-interpret "242.53" as a
# number and silently change the cell's stored type.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# --- simple price updates (column D only) ---
Set-TextValue "D2"  "242.53"
Set-TextValue "D4"  "5.414"
Set-TextValue "D5"  "0.05888"
Set-TextValue "D6"  "3.435"
Set-TextValue "D7"  "6.528"
Set-TextValue "D9"  "0.9288"
Set-TextValue "D10" "0.1417"
Set-TextValue "D11" "0.07382"
Set-TextValue "D12" "0.03290"
Set-TextValue "D13" "0.03069"
Set-TextValue "D14" "0.09356"
Set-TextValue "D15" "3.846"
Set-TextValue "D16" "0.001587"
Set-TextValue "D17" "0.04679"

# --- rows 18-24: coins reshuffled (each row's data shifted), with a few
#     price/volume tweaks layered on top ---
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005892"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.005953"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.001259"
$ws.Range("E20").Value = "19BitKanKANBestin24h"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D21" "0.004900"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.00006802"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.563"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.141"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# --- more simple price updates further down the sheet ---
Set-TextValue "D40" "0.03975"
Set-TextValue "D41" "0.006180"
Set-TextValue "D43" "0.002571"
Set-TextValue "D44" "0.008747"
Set-TextValue "D45" "0.00005189"
Set-TextValue "D47" "0.6702"
Set-TextValue "D48" "0.002337"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"
